$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 22-30: regcntr_id continues 10002..10010, machine_id continues 10021..10029
$regcntrStart = 10002
$machineStart = 10021

for ($i = 0; $i -lt 9; $i++) {
    $r = 22 + $i
    $ws.Cells.Item($r, 1).Value = $regcntrStart + $i
    $ws.Cells.Item($r, 2).Value = $machineStart + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Selection state matching the saved view
$ws.Range("B22:B30").Select()

# Page setup (portrait orientation, as saved)
$ws.PageSetup.Orientation = 1

$wb.Save()
